$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 251 (this shifts existing rows 251-272 down to 252-273)
$ws.Rows.Item(251).Insert()

# Populate the newly inserted row 251 with the new record
$ws.Range("A251").Value = 5
$ws.Range("B251").Value = "Macroferia Regional de Talca"
$ws.Range("C251").Value = "Maule"
$ws.Range("D251").Value = 45194
$ws.Range("E251").Value = 7
$ws.Range("F251").Value = 100112031
$ws.Range("G251").Value = "Poroto verde"
$ws.Range("H251").Value = "Sin especificar"
$ws.Range("I251").Value = "Primera"
$ws.Range("J251").Value = 150
$ws.Range("K251").Value = 30000
$ws.Range("L251").Value = 30000
$ws.Range("M251").Value = 30000
$ws.Range("N251").Value = "$/malla 25 kilos"
$ws.Range("O251").Value = "Perú"
$ws.Range("P251").Value = 1200
$ws.Range("Q251").Value = 25
$ws.Range("R251").Value = "Hortaliza"

# Ensure date cell keeps the same date number format used by the other rows in column D
$ws.Range("D251").NumberFormat = $ws.Range("D252").NumberFormat
